$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain value changes (no formulas involved) ---
$ws.Range("F12").Value = 1305498328.7699957
$ws.Range("G12").Value = 1240524717

$ws.Range("F13").Value = 325268233.58999997
$ws.Range("G13").Value = 3198194883.1999998

$ws.Range("G14").Value = 34063116.799999997

$ws.Range("F16").Value = -53616441.74000001

# --- F18 becomes a formula (was a plain value before) ---
$ws.Range("F18").Formula = "=SUM(F12:F17)"
# G18 already holds =SUM(G12:G17); it will recompute automatically.

$ws.Range("F19").Value = -412700000
$ws.Range("G19").Value = -379300000

# --- F21 becomes a formula (was a plain value before) ---
$ws.Range("F21").Formula = "=SUM(F18:F20)"
# G21 already holds =SUM(G18:G20); it will recompute automatically.

# F23 and G23 already hold SUM formulas; they recompute automatically.
# F25 and G25 already hold division formulas; they recompute automatically.

$ws.Range("F26").Value = 1026703455.3810816
$ws.Range("G26").Value = 1029174575

# F28, G28, and C29 already hold formulas; they recompute automatically.

$excel.Calculate()
